$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.2874596666666667
$ws.Range("H2").Value = 0.862379
$ws.Range("I2").Value = 0.113106005361649
$ws.Range("J2").Value = 0.113106005361649
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.2585373333333333
$ws.Range("N2").Value = 0.775612
$ws.Range("O2").Value = 0.2883652589930572
$ws.Range("P2").Value = 0.2883652589930571
$ws.Range("Q2").Value = 0.0743190556608889
$ws.Range("R2").Value = 0.668871500948
$ws.Range("S2").Value = 0.03261584252978202
$ws.Range("T2").Value = 0.03261584252978202

# Row 3
$ws.Range("G3").Value = 0.2874596666666667
$ws.Range("H3").Value = 0.862379
$ws.Range("I3").Value = 0.113106005361649
$ws.Range("J3").Value = 0.113106005361649
$ws.Range("O3").Value = 0.1945740134722046
$ws.Range("P3").Value = 0.1945740134722046
$ws.Range("Q3").Value = 0.05014666811077778
$ws.Range("R3").Value = 0.451320012997
$ws.Range("S3").Value = 0.02200748941102473
$ws.Range("T3").Value = 0.02200748941102473

# Row 4
$ws.Range("D4").Value = "MuSCs"
$ws.Range("G4").Value = 0.2874596666666667
$ws.Range("H4").Value = 0.862379
$ws.Range("I4").Value = 0.113106005361649
$ws.Range("J4").Value = 0.113106005361649
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.032643
$ws.Range("N4").Value = 0.097929
$ws.Range("O4").Value = 0.03640908269589833
$ws.Range("P4").Value = 0.03640908269589833
$ws.Range("Q4").Value = 0.009383545899
$ws.Range("R4").Value = 0.084451913091
$ws.Range("S4").Value = 0.004118085902614998
$ws.Range("T4").Value = 0.004118085902614999

# Row 5
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("G5").Value = 0.2874596666666667
$ws.Range("H5").Value = 0.862379
$ws.Range("I5").Value = 0.113106005361649
$ws.Range("J5").Value = 0.113106005361649
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.430934
$ws.Range("N5").Value = 1.292802
$ws.Range("O5").Value = 0.4806516448388399
$ws.Range("P5").Value = 0.4806516448388399
$ws.Range("Q5").Value = 0.1238761439953333
$ws.Range("R5").Value = 1.114885295958
$ws.Range("S5").Value = 0.05436458751822723
$ws.Range("T5").Value = 0.05436458751822724

# Row 6
$ws.Range("G6").Value = 2.225929666666667
$ws.Range("H6").Value = 6.677789
$ws.Range("I6").Value = 0.8758307408204057
$ws.Range("J6").Value = 0.8758307408204058
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.2585373333333333
$ws.Range("N6").Value = 0.775612
$ws.Range("O6").Value = 0.2883652589930572
$ws.Range("P6").Value = 0.2883652589930571
$ws.Range("Q6").Value = 0.5754859202075555
$ws.Range("R6").Value = 5.179373281867999
$ws.Range("S6").Value = 0.2525591584107574
$ws.Range("T6").Value = 0.2525591584107574

# Row 7
$ws.Range("G7").Value = 2.225929666666667
$ws.Range("H7").Value = 6.677789
$ws.Range("I7").Value = 0.8758307408204057
$ws.Range("J7").Value = 0.8758307408204058
$ws.Range("O7").Value = 0.1945740134722046
$ws.Range("P7").Value = 0.1945740134722046
$ws.Range("Q7").Value = 0.3883082365141111
$ws.Range("R7").Value = 3.494774128627
$ws.Range("S7").Value = 0.1704139023637605
$ws.Range("T7").Value = 0.1704139023637606

# Row 8
$ws.Range("D8").Value = "MuSCs"
$ws.Range("G8").Value = 2.225929666666667
$ws.Range("H8").Value = 6.677789
$ws.Range("I8").Value = 0.8758307408204057
$ws.Range("J8").Value = 0.8758307408204058
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.032643
$ws.Range("N8").Value = 0.097929
$ws.Range("O8").Value = 0.03640908269589833
$ws.Range("P8").Value = 0.03640908269589833
$ws.Range("Q8").Value = 0.072661022109
$ws.Range("R8").Value = 0.653949198981
$ws.Range("S8").Value = 0.03188819387014005
$ws.Range("T8").Value = 0.03188819387014005

# Row 9
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("G9").Value = 2.225929666666667
$ws.Range("H9").Value = 6.677789
$ws.Range("I9").Value = 0.8758307408204057
$ws.Range("J9").Value = 0.8758307408204058
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.430934
$ws.Range("N9").Value = 1.292802
$ws.Range("O9").Value = 0.4806516448388399
$ws.Range("P9").Value = 0.4806516448388399
$ws.Range("Q9").Value = 0.9592287749753332
$ws.Range("R9").Value = 8.633058974777999
$ws.Range("S9").Value = 0.4209694861757477
$ws.Range("T9").Value = 0.4209694861757478

# Row 10
$ws.Range("I10").Value = 0.01106325381794526
$ws.Range("J10").Value = 0.01106325381794526
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.2585373333333333
$ws.Range("N10").Value = 0.775612
$ws.Range("O10").Value = 0.2883652589930572
$ws.Range("P10").Value = 0.2883652589930571
$ws.Range("Q10").Value = 0.007269380380444444
$ws.Range("R10").Value = 0.065424423424
$ws.Range("S10").Value = 0.003190258052517713
$ws.Range("T10").Value = 0.003190258052517713

# Row 11
$ws.Range("I11").Value = 0.01106325381794526
$ws.Range("J11").Value = 0.01106325381794526
$ws.Range("O11").Value = 0.1945740134722046
$ws.Range("P11").Value = 0.1945740134722046
$ws.Range("S11").Value = 0.002152621697419299
$ws.Range("T11").Value = 0.0021526216974193

# Row 12
$ws.Range("D12").Value = "MuSCs"
$ws.Range("I12").Value = 0.01106325381794526
$ws.Range("J12").Value = 0.01106325381794526
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.032643
$ws.Range("N12").Value = 0.097929
$ws.Range("O12").Value = 0.03640908269589833
$ws.Range("P12").Value = 0.03640908269589833
$ws.Range("Q12").Value = 0.0009178341119999998
$ws.Range("R12").Value = 0.008260507008
$ws.Range("S12").Value = 0.0004028029231432819
$ws.Range("T12").Value = 0.000402802923143282

# Row 13
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("I13").Value = 0.01106325381794526
$ws.Range("J13").Value = 0.01106325381794526
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.430934
$ws.Range("N13").Value = 1.292802
$ws.Range("O13").Value = 0.4806516448388399
$ws.Range("P13").Value = 0.4806516448388399
$ws.Range("Q13").Value = 0.01211671492266667
$ws.Range("R13").Value = 0.109050434304
$ws.Range("S13").Value = 0.005317571144864965
$ws.Range("T13").Value = 0.005317571144864965

